$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "question" header/column to "statement"
$ws.Range("B1").Value = "statement"

# Update the selection to match the target state
$ws.Range("B1").Select()
